$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: copy the date-cell formatting from A4 so A5 keeps the same style
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A5").Value = 41426
$ws.Range("B5").Value = 0.75
$ws.Range("D5").Value = "Implementation design of new sync objects "

# Match the diff's updated selection (active cell now D5)
$ws.Range("D5").Select()
